# Updates Price (column D) and Volume(1h) (column E) cells on Sheet1
# for rows 2-51, per the latest cryptos-list scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as literal text even though
# it looks like a plain number (e.g. "519.46"), so Excel does not silently
# convert it to a numeric cell. Number format is restored to the default
# afterwards so no visible/style change is introduced.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "60.160.61"
$ws.Range("E2").Value = "  -1.23%  "
# Row 3
$ws.Range("D3").Value = "2.627.13"
$ws.Range("E3").Value = "  +1.04%  "
# Row 4
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
Set-TextValue $ws.Range("D5") "519.46"
$ws.Range("E5").Value = "  -0.56%  "
# Row 6
Set-TextValue $ws.Range("D6") "147.94"
$ws.Range("E6").Value = "  -4.39%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  +0.03%  "
# Row 8
$ws.Range("E8").Value = "  -3.53%  "
# Row 9
$ws.Range("D9").Value = "2.632.49"
$ws.Range("E9").Value = "  +0.97%  "
# Row 10
$ws.Range("E10").Value = "  -5.58%  "
# Row 11
$ws.Range("E11").Value = "  -0.51%  "
# Row 12
$ws.Range("E12").Value = "  -2.17%  "
# Row 13
$ws.Range("E13").Value = "  -0.70%  "
# Row 14
$ws.Range("D14").Value = "3.083.52"
$ws.Range("E14").Value = "  +0.93%  "
# Row 15
$ws.Range("D15").Value = "60.145.37"
$ws.Range("E15").Value = "  -1.28%  "
# Row 16
Set-TextValue $ws.Range("D16") "21.18"
$ws.Range("E16").Value = "  -2.44%  "
# Row 17
$ws.Range("E17").Value = "  -2.02%  "
# Row 18
$ws.Range("D18").Value = "2.629.74"
$ws.Range("E18").Value = "  +0.98%  "
# Row 19
Set-TextValue $ws.Range("D19") "4.63"
$ws.Range("E19").Value = "  -2.16%  "
# Row 20
Set-TextValue $ws.Range("D20") "340.34"
$ws.Range("E20").Value = "  -3.52%  "
# Row 21
Set-TextValue $ws.Range("D21") "10.41"
$ws.Range("E21").Value = "  -1.41%  "
# Row 22
Set-TextValue $ws.Range("D22") "6.12"
$ws.Range("E22").Value = "  -1.56%  "
# Row 23
$ws.Range("E23").Value = "  -0.37%  "
# Row 24
Set-TextValue $ws.Range("D24") "61.17"
$ws.Range("E24").Value = "  +0.13%  "
# Row 25
Set-TextValue $ws.Range("D25") "0.418"
$ws.Range("E25").Value = "  -2.03%  "
# Row 26
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  -0.08%  "
# Row 27
$ws.Range("E27").Value = "  -3.69%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0809"
$ws.Range("E28").Value = "  -4.63%  "
# Row 29
Set-TextValue $ws.Range("D29") "7.04"
$ws.Range("E29").Value = "  -4.42%  "
# Row 30
$ws.Range("E30").Value = "  -0.05%  "
# Row 31
$ws.Range("E31").Value = "  -1.05%  "
# Row 32
Set-TextValue $ws.Range("D32") "5.97"
$ws.Range("E32").Value = "  -5.57%  "
# Row 33
Set-TextValue $ws.Range("D33") "18.91"
$ws.Range("E33").Value = "  -2.34%  "
# Row 34
Set-TextValue $ws.Range("D34") "149.93"
$ws.Range("E34").Value = "  +0.58%  "
# Row 35
$ws.Range("E35").Value = "  -6.82%  "
# Row 36
Set-TextValue $ws.Range("D36") "0.920"
$ws.Range("E36").Value = "  -2.32%  "
# Row 37
Set-TextValue $ws.Range("D37") "1.14"
$ws.Range("E37").Value = "  -5.67%  "
# Row 38
Set-TextValue $ws.Range("D38") "0.858"
$ws.Range("E38").Value = "  +1.18%  "
# Row 39
Set-TextValue $ws.Range("D39") "36.69"
$ws.Range("E39").Value = "  +0.62%  "
# Row 40
$ws.Range("E40").Value = "  -4.70%  "
# Row 41
$ws.Range("E41").Value = "  -4.06%  "
# Row 42
Set-TextValue $ws.Range("D42") "290.75"
$ws.Range("E42").Value = "  +1.14%  "
# Row 43
Set-TextValue $ws.Range("D43") "0.628"
$ws.Range("E43").Value = "  +0.55%  "
# Row 44
Set-TextValue $ws.Range("D44") "0.0999"
$ws.Range("E44").Value = "  -1.22%  "
# Row 45
$ws.Range("E45").Value = "  +0.15%  "
# Row 47
Set-TextValue $ws.Range("D47") "19.41"
$ws.Range("E47").Value = "  -0.80%  "
# Row 48
$ws.Range("E48").Value = "  +0.82%  "
# Row 49
$ws.Range("E49").Value = "  -2.24%  "
# Row 50
Set-TextValue $ws.Range("D50") "4.65"
$ws.Range("E50").Value = "  -4.63%  "
# Row 51
$ws.Range("D51").Value = "1.957.99"
$ws.Range("E51").Value = "  -0.11%  "
